# TimeLog_(Kubra).xlsx - "Add files via upload"
#
# The author filled in four more weekly timesheet rows (project weeks 8-11,
# rows 14-17 on Tabelle1) that had previously been left blank: hours worked
# (col E) and the activity description (col F). Filling these cells also
# appends their four distinct description strings to the shared-string
# table.
#
# NOTE on ordering: the shared-string table in the target file has
# "Meeting with other team mates" (row 17) inserted *before*
# "Setup of all environments on the new Laptop" (row 16), so we enter the
# data in that same order (row 17 before row 16) to reproduce the exact
# shared-string indices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Project week 8 (row 14)
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = "Learning SpringBoot Framework"

# Project week 9 (row 15)
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = "Stored Procedure for Panic Button"

# Project week 11 (row 17) - entered before row 16, see note above
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = "Meeting with other team mates"

# Project week 10 (row 16)
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = "Setup of all environments on the new Laptop"

# Reflect the author's final cursor position/selection in the sheet.
$ws.Range("F18").Select() | Out-Null
